$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4999.5
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3752.8918
$ws.Range("I32").Value = 3752.8918
$ws.Range("K32").Value = 3752.8918
$ws.Range("M32").Value = -3465.8918

$ws.Range("H45").Value = 3099.8
$ws.Range("J45").Value = 3499.75
$ws.Range("L45").Value = 3499.75
$ws.Range("N45").Value = -4253.75

$ws.Range("H61").Value = 3060.7334
$ws.Range("I61").Value = 1489
$ws.Range("J61").Value = 4857
$ws.Range("K61").Value = 1489
$ws.Range("L61").Value = 4857
$ws.Range("M61").Value = -1277
$ws.Range("N61").Value = -5281

$ws.Range("H63").Value = 10753
$ws.Range("I63").Value = 1500
$ws.Range("K63").Value = 1500
$ws.Range("M63").Value = -814

$ws.Range("H66").Value = 10753
$ws.Range("I66").Value = 1500
$ws.Range("K66").Value = 7500
$ws.Range("M66").Value = -4068

$ws.Range("H74").Value = 1186.8182
$ws.Range("I74").Value = 1247
$ws.Range("J74").Value = 585
$ws.Range("K74").Value = 1247
$ws.Range("L74").Value = 585
$ws.Range("M74").Value = -373
$ws.Range("N74").Value = -2333

$ws.Range("H77").Value = 1186.8182
$ws.Range("I77").Value = 1247
$ws.Range("J77").Value = 585
$ws.Range("K77").Value = 6235
$ws.Range("L77").Value = 2925
$ws.Range("M77").Value = -1867
$ws.Range("N77").Value = -11661

$ws.Range("H136").Value = 3060.7334
$ws.Range("I136").Value = 1489
$ws.Range("J136").Value = 4857
$ws.Range("K136").Value = 4467
$ws.Range("L136").Value = 14571
$ws.Range("M136").Value = -1917
$ws.Range("N136").Value = -19671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 796.3333
$ws.Range("I107").Value = 797
$ws.Range("K107").Value = 797
$ws.Range("M107").Value = 1123

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3609.762
$ws.Range("I31").Value = 1216
$ws.Range("J31").Value = 4806.643
$ws.Range("K31").Value = 1216
$ws.Range("L31").Value = 4806.643
$ws.Range("M31").Value = -921
$ws.Range("N31").Value = -5396.643

$ws.Range("H34").Value = 3609.762
$ws.Range("I34").Value = 1216
$ws.Range("J34").Value = 4806.643
$ws.Range("K34").Value = 1216
$ws.Range("L34").Value = 4806.643
$ws.Range("M34").Value = -1014
$ws.Range("N34").Value = -5210.643

$ws.Range("H58").Value = 3000
$ws.Range("I58").Value = 3000
$ws.Range("K58").Value = 3000
$ws.Range("M58").Value = -2797

$ws.Range("H86").Value = 8640.727999999999
$ws.Range("I86").Value = 6262.3335
$ws.Range("K86").Value = 6262.3335
$ws.Range("M86").Value = -5139.3335

$ws.Range("H89").Value = 8640.727999999999
$ws.Range("I89").Value = 6262.3335
$ws.Range("K89").Value = 31311.6675
$ws.Range("M89").Value = -25695.6675

$ws.Range("H94").Value = 1336.75
$ws.Range("I94").Value = 1281
$ws.Range("K94").Value = 1281
$ws.Range("M94").Value = -830

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -6450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1184.3334
$ws.Range("I5").Value = 1221.2
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 3663.6
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -3551.6
$ws.Range("N5").Value = -3224

$ws.Range("H109").Value = 4193.3477
$ws.Range("I109").Value = 1289.4
$ws.Range("K109").Value = 3868.2
$ws.Range("M109").Value = -2828.2

$ws.Range("H135").Value = 1184.3334
$ws.Range("I135").Value = 1221.2
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 10990.8
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -8455.800000000001
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6179.4287
$ws.Range("J122").Value = 1200
$ws.Range("L122").Value = 3600
$ws.Range("N122").Value = -8500

$ws.Range("H132").Value = 4066
$ws.Range("I132").Value = 2200
$ws.Range("K132").Value = 6600
$ws.Range("M132").Value = -4070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 2003
$ws.Range("I19").Value = 2003
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2003
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1833
$ws.Range("N19").ClearContents()

$ws.Range("H46").Value = 2359.3
$ws.Range("I46").Value = 1598.8334
$ws.Range("K46").Value = 1598.8334
$ws.Range("M46").Value = -1410.8334

$ws.Range("H64").Value = 20150
$ws.Range("J64").Value = 20150
$ws.Range("L64").Value = 20150
$ws.Range("N64").Value = -20600

$ws.Range("H67").Value = 20150
$ws.Range("J67").Value = 20150
$ws.Range("L67").Value = 20150
$ws.Range("N67").Value = -21710

$ws.Range("H68").Value = 2000
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1251
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 2000
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -6256
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1448.1666
$ws.Range("I96").Value = 1622.5
$ws.Range("J96").Value = 1099.5
$ws.Range("K96").Value = 1622.5
$ws.Range("L96").Value = 1099.5
$ws.Range("M96").Value = -249.5
$ws.Range("N96").Value = -3845.5
